$d = $word.ActiveDocument

# Locate the last step of "modo de preparo"
# ("6= Asse em forno preaquecido 180C por cerca de 40 minutos.") using
# Find so the script does not depend on paragraph indices.
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "6= Asse em forno preaquecido 180",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the end of the 'modo de preparo' section."
}

# Insertion point: right after the matched text (i.e. right before that
# paragraph's own paragraph mark), so the new content is injected as new
# paragraphs right after it and right before the pre-existing trailing
# empty paragraph.
$insertionPoint = $d.Range($searchRange.End, $searchRange.End)

$newContentXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t xml:space="preserve">Cobertura: </w:t>
            </w:r>
            <w:r>
              <w:t xml:space="preserve"> Chocolate meio amargo.</w:t>
            </w:r>
          </w:p>
          <w:p/>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:b/>
                <w:bCs/>
              </w:rPr>
              <w:t xml:space="preserve">DICAS DE PREPARO </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t xml:space="preserve">- Para um bolo mais fofo, peneire a farinha de trigo. </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>- Você pode adicionar nozes picadas à mesa para dar um toque especial.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($newContentXml) | Out-Null
